# Update the "want to go" count (想去人数, column F) for two rows on both the
# "展览" and "全部类型" worksheets, matching the upstream data refresh:
#   F2: 309 -> 310
#   F5: 274 -> 275

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 310
    $ws.Range("F5").Value = 275
}
